$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "08.00 AM"
$ws.Range("D2").Value = 45092
$ws.Range("H10").Select()
